$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new row of data ("數" / "v" / "物品") is inserted above the old row 8,
# shifting the previous rows 8-11 down to rows 9-12.
$ws.Range("A8").Value = "數"
$ws.Range("B8").Value = "v"
$ws.Range("C8").Value = "物品"

$ws.Range("A9").Value = "多"
$ws.Range("B9").Value = "bigger"
$ws.Range("C9").Value = "大"

$ws.Range("B10").Value = "int"
$ws.Range("C10").Value = "數量"

# "5" looks numeric, so a plain .Value assignment would store it as a
# number instead of text (t="s"). Round-trip it through a text formula
# and Paste Values so it lands back as a shared string, matching the
# original cell's "t=s" storage.
$a10 = $ws.Range("A10")
$a10.Formula = "=""5"""
$a10.Copy()
$a10.PasteSpecial(-4163)
$excel.CutCopyMode = $false

$ws.Range("A11").Value = "。"
# B11 already holds "x" and is untouched by the diff.

# A12 already holds the newline string and is untouched by the diff;
# only B12 is a genuinely new cell.
$ws.Range("B12").Value = "x"

# Reflect the new selection recorded in the sheet view (best achievable
# approximation of the multi-area selection, keyed on its active cell).
$ws.Range("G1").Select()
